$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 26 (shifts existing rows 26-105 down to 27-106)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new daily price entry
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44575
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 100112022
$ws.Range("G26").Value = "Arveja Verde"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = 30000
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región de La Araucanía"
$ws.Range("P26").Value = 1200
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
